# Taakstructuur.xlsx - "Vergaderverslagen klaar + taakstructuur aangepast"
#
# Updates a handful of task-status cells (column E) in the CAD-model
# section, removes the two blank "group separator" rows that used to sit
# above the RAPPORTERING/VERSLAG (row 51) and LABVIEW (row 101) group
# headers, and restores the previously-selected cell/viewport.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Status column (E) updates: CAD MODEL sub-tasks (2.1.x) ---------------
# 'bezig' (in progress) -> resolved to either 'OK' or 'niet OK'
$ws.Range("E16").Value = "OK"
$ws.Range("E17").Value = "OK"
$ws.Range("E18").Value = "niet OK"
$ws.Range("E19").Value = "niet OK"
$ws.Range("E20").Value = "OK"
$ws.Range("E21").Value = "niet OK"
$ws.Range("E22").Value = "OK"
$ws.Range("E23").Value = "niet OK"
$ws.Range("E24").Value = "niet OK"
$ws.Range("E25").Value = "OK"

# Group headers for RAPPORTERING/VERSLAG (4) and LABVIEW (6) moved back to
# 'bezig' (work resumed on those sections)
$ws.Range("E31").Value = "bezig"
$ws.Range("E47").Value = "bezig"

# --- Row 50: drop the blank separator row's border formatting, leaving ---
# --- only a plain, centre-aligned blank cell in B50 (matches the style  --
# --- used by the other blank spacer rows, e.g. row 29/98/99).          ---
$ws.Range("B50:E50").Clear()
$ws.Range("B50").HorizontalAlignment = -4108   ' xlCenter

# --- Row 100: this separator row is removed entirely -----------------------
$ws.Range("B100:E100").Clear()

# --- Restore the active selection/viewport used in the saved session ------
$ws.Range("E48").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 73
$win.ScrollColumn = 1
